$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest scraped quotes, row by row, as published by the
# "Updated cryptos list" GitHub Actions job.
#
# Two price cells (D32, D45) are written with a leading apostrophe so Excel
# keeps the significant trailing zero ("1.10", "19.30") as literal text
# instead of silently normalising them to the numbers 1.1 / 19.3.

$ws.Range("D2").Value = '59.760.08'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '2.364.81'
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '557.05'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").Value = '133.04'
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = '5.62'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '0.341'
$ws.Range("E12").Value = '  -3.53%  '
$ws.Range("D13").Value = '24.19'
$ws.Range("E13").Value = '  -4.23%  '
$ws.Range("D14").Value = '2.782.82'
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("D15").Value = '59.701.79'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '2.370.52'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '11.05'
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").Value = '4.45'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").Value = '319.71'
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").Value = '6.63'
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '64.18'
$ws.Range("E23").Value = '  -2.64%  '
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '8.36'
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").Value = '0.0₃0755'
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("D30").Value = '170.72'
$ws.Range("E30").Value = '  +0.84%  '
$ws.Range("D31").Value = '6.06'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = '  +9.44%  '
$ws.Range("D33").Value = '0.397'
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").Value = '18.08'
$ws.Range("E34").Value = '  -2.85%  '
$ws.Range("D36").Value = '1.31'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '4.11'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("D40").Value = '316.99'
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("D41").Value = '38.58'
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("D42").Value = '144.35'
$ws.Range("E42").Value = '  +3.06%  '
$ws.Range("D43").Value = '3.52'
$ws.Range("E43").Value = '  -3.73%  '
$ws.Range("D44").Value = '0.0964'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = "'19.30"
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("D46").Value = '0.0509'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").Value = '0.567'
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").Value = '1.53'
$ws.Range("E51").Value = '  -1.93%  '
